$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update data row 2
$ws.Range("B2").Value = 42.79353601253937
$ws.Range("C2").Value = -12.452003414729768
$ws.Range("D2").Value = -5.454241636103033
$ws.Range("E2").Value = 1.7946541736515533

# Update data row 3
$ws.Range("B3").Value = 40.635045313971894
$ws.Range("C3").Value = 7.226494459758567
$ws.Range("D3").Value = -23.671649939920549
$ws.Range("E3").Value = 15.667814657011764

# Update the selection to reflect the new used range selection (B1:E3)
$ws.Range("B1:E3").Select()
